$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Change 1: "Image difference is used for selective " -> split off the
# word "selective" into its own bold run (matching the bold style used
# for "spectral enhancement" / "change detection" / "removal of
# background illumination bias" later in the same sentence), leaving a
# trailing space in a separate (non-bold) run.
# -------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Image difference is used for selective", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $end = $rng.End
    $wordStart = $end - 9   # length of "selective"
    $sub = $d.Range($wordStart, $end)
    $sub.Bold = 1
    $sub.BoldBi = 1
}

# -------------------------------------------------------------------------
# Change 2: merge the two runs
#   "For NDVI and all normalised difference indices, the key part of the
#    recipe is that it is a two-band difference image normalised by the
#    sum of the same two bands" + "."
# into a single run (same text, same formatting).
# -------------------------------------------------------------------------
$rng2 = $d.Content
$old2 = "For NDVI and all normalised difference indices, the key part of the recipe is that it is a two-band difference image normalised by the sum of the same two bands."
$rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# -------------------------------------------------------------------------
# Change 3: merge the six runs making up the "For the mineral indices'
# formula..." sentence into a single run (same text, same formatting).
# -------------------------------------------------------------------------
$rng3 = $d.Content
$rsquo = [char]0x2019
$old3 = "For the mineral indices" + $rsquo + " formula, the key thing here is the subtraction of the minimum DN value in that band from every individual DN value in that band. This subtraction forces the histogram value range of both bands to start at zero. So, the division between them is more representative."
$rng3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2) | Out-Null
